$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Country" column before the current column L ("State") ---
# This shifts the old L:S (State .. Primary Delegate Phone Type) to M:T.
$ws.Columns("L:L").Insert(-4161)

# New L1/L2 inherited K's style (s=5 / s=3). Header needs the "s=4" look (same as A1),
# so copy formatting from A1 onto L1. L2 already matches the blank/body style (s=3, like A2).
$ws.Range("A1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("L1").Value = "Country"
$ws.Range("L2").Value = ""

# --- Append two new trailing columns: U (Primary Delegate Year of Birth) and
# V (Primary Delegate is over the age of legal majority) ---

# Header style s=5 (like B1); body style s=3 (like A2).
$ws.Range("B1").Copy()
$ws.Range("U1").PasteSpecial(-4122)
$ws.Range("V1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Copy()
$ws.Range("U2").PasteSpecial(-4122)
$ws.Range("V2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("U1").Value = "Primary Delegate Year of Birth"
$ws.Range("U2").Value = "YYYY (formatted as text)"

$ws.Range("V1").Value = "Primary Delegate is over the age of legal majority"
$ws.Range("V2").Value = "Enter one of these values:         Yes     `n No"

# --- Row heights / selection / dimension tracking ---
$ws.Rows("1:1").RowHeight = 78
$ws.Rows("2:2").RowHeight = 124.8

$ws.Range("L2").Select()
